# "#Update plan with table"
#
# 1) Footer/date placeholder on the slide master + every slide layout:
#    "1/14/2019" -> "01/14/2019" (zero-padded month).
# 2) On slide 15 ("Update plan with table"), nudge five of the table-header
#    rectangles (Agency / Last name / Invoice number columns) a few EMU to
#    better align with the new table.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "1/14/2019") {
                $shp.TextFrame.TextRange.Text = "01/14/2019"
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 15: reposition the five header rectangles over the new table.
$moves = @{
    "Rectangle 10" = @{ Left = 264.94236755371094; Top = 97.42133712768555 }
    "Rectangle 11" = @{ Left = 263.74250793457037; Top = 156.7223587036133 }
    "Rectangle 12" = @{ Left = 265.17173767089844; Top = 196.20984649658206 }
    "Rectangle 21" = @{ Left = 383.47187805175787; Top = 155.91118621826172 }
    "Rectangle 23" = @{ Left = 383.23432922363287; Top = 198.35717010498047 }
}

$slide = $p.Slides.Item(15)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($moves.ContainsKey($shp.Name)) {
        $shp.Left = $moves[$shp.Name].Left
        $shp.Top = $moves[$shp.Name].Top
    }
}
